# Finished labeling LLM_corrections_check v2.0
# Applies the remaining row 53-101 misinterpretation labels/comments,
# highlights a handful of tricky/flagged rows in yellow, turns on
# AutoFilter for the data range, and resets the sheet view back to A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("llm_correction_checked")

# ---------------------------------------------------------------------
# 1) Column D ("misinterpretation?" flag) for rows 53-101.
#    Default every still-blank row to 0, then override the rows that
#    were coded as 1 (actual misinterpretations).
# ---------------------------------------------------------------------
for ($r = 53; $r -le 101; $r++) {
    $ws.Cells.Item($r, 4).Value = 0
}

$misinterpretationRows = @(75, 78, 80, 81, 101)
foreach ($r in $misinterpretationRows) {
    $ws.Cells.Item($r, 4).Value = 1
}

# ---------------------------------------------------------------------
# 2) Column E ("comment") new notes. Order matters: it controls the
#    order new entries are appended to the shared-string table.
# ---------------------------------------------------------------------
$ws.Range("E53").Value = 'borderline: "not statistically different"; I allow ''statistically'' as a synonym for ''significantly'' here, knowing that it''s not ideal'
$ws.Range("E63").Value = "a little strange, but still correct"
$ws.Range("E68").Value = "reporting correlations like descriptives, ignoring the significance"
$ws.Range("E75").Value = "very tricky; sure, they're compatible with no effect, but also a range of non-zero effects; I count this as a misinterpretation"
$ws.Range("E78").Value = "very tricky; sure, they're compatible with no effect, but also a range of non-zero effects; I count this as a misinterpretation"
$ws.Range("E80").Value = "Check code again, little strange"
$ws.Range("E101").Value = "Again, it might be consistent with a zero effect, but also with all different kinds of other effects"
$ws.Range("E81").Value = "No misinterpretation as no effect, but the LLM removed the p values and reported the effects as descriptives instead? Without the CIs it sounds like the effects were significant…"

# ---------------------------------------------------------------------
# 3) Highlight the flagged / borderline comments in yellow.
# ---------------------------------------------------------------------
$highlightRows = @(24, 41, 75, 78, 80, 81, 101)
foreach ($r in $highlightRows) {
    $ws.Cells.Item($r, 5).Interior.Color = 65535
}

# ---------------------------------------------------------------------
# 4) Row 81 grew to three lines of wrapped text -> taller row.
# ---------------------------------------------------------------------
$ws.Range("A81").EntireRow.RowHeight = 75

# ---------------------------------------------------------------------
# 5) Turn on AutoFilter for the whole table and register the
#    corresponding hidden workbook-level filter-database name.
# ---------------------------------------------------------------------
[void]$ws.Range("A1:E101").AutoFilter()
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=llm_correction_checked!`$A`$1:`$E`$101")
$filterName.Visible = $false

# ---------------------------------------------------------------------
# 6) Scroll/selection back to the top of the sheet.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("A1").Select()

Write-Output "edit complete"
